$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 9: L9 gains an "x"
$ws.Range("L9").Value = "x"

# Row 10: several cells change from "x" to "s", and two "x" cells are cleared,
# while an empty cell gains an "s"
$ws.Range("E10").Value = "s"
$ws.Range("F10").Value = "s"
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = "s"
$ws.Range("K10").Value = "s"

# Update the active selection to K5
$ws.Range("K5").Select() | Out-Null
